# Rotate the data rows 7-19 on the Artfynd sheet.
# Each destination row ends up with the full contents (all columns A:AY)
# that a particular source row had before the edit. This reproduces the
# "automatic update" diff, where every field of a record (id, coordinates,
# species names, dates, etc.) moved together as a unit to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (values captured from the *original* sheet)
$mapping = @{
    7  = 8
    8  = 9
    9  = 11
    10 = 12
    11 = 13
    12 = 14
    13 = 15
    14 = 16
    15 = 17
    16 = 18
    17 = 19
    18 = 7
    19 = 10
}

$firstCol = "A"
$lastCol  = "AY"

# 1) Snapshot every source row's full row of values BEFORE any writes happen,
#    since several rows are both a source and a destination (it's one big
#    13-row permutation cycle).
$snapshots = @{}
foreach ($row in 7..19) {
    $rng = $ws.Range("$firstCol$row" + ":" + "$lastCol$row")
    $snapshots[$row] = $rng.Value2
}

# 2) Write each destination row from its recorded source snapshot.
foreach ($destRow in 7..19) {
    $srcRow = $mapping[$destRow]
    $destRng = $ws.Range("$firstCol$destRow" + ":" + "$lastCol$destRow")
    $destRng.Value = $snapshots[$srcRow]
}

# 3) Columns Y ("Startdatum") and AA ("Slutdatum") hold plain date-looking
#    text (e.g. "2021-10-21") in the source workbook, not real dates. The
#    bulk array write above lets Excel "smart convert" those strings into
#    date serials, so re-assert them as literal text afterwards, using the
#    values captured in the original (pre-write) snapshot.
$colIndex = @{ "Y" = 25; "AA" = 27 }
foreach ($destRow in 7..19) {
    $srcRow = $mapping[$destRow]
    foreach ($col in @("Y", "AA")) {
        $cell = $ws.Range("$col$destRow")
        $cell.NumberFormat = "@"
        $cell.Value = $snapshots[$srcRow][1, $colIndex[$col]]
        $cell.NumberFormat = "General"
    }
}
